$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (column A values are text, like the existing row 21,
# so format as text before assigning and then clear the formatting override
# so the cell keeps the default style while remaining text).
$newRows = @(
    @{ Row = 22; A = "20072900"; B = 0;             C = 0;             D = 0 },
    @{ Row = 23; A = "20073000"; B = 165000000000;  C = 186000000000;  D = -21000000000 },
    @{ Row = 24; A = "20073100"; B = 196000000000;  C = 151000000000;  D = 45000000000 },
    @{ Row = 25; A = "20080100"; B = 0;             C = 0;             D = 0 },
    @{ Row = 26; A = "20080200"; B = 0;             C = 0;             D = 0 },
    @{ Row = 27; A = "20080300"; B = 0;             C = 0;             D = 0 }
)

foreach ($r in $newRows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.A
    $cellA.ClearFormats()

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
